$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 126, shifting existing rows 126-240 down to 128-242.
$ws.Rows("126:127").Insert()

# Row 126 (new): Primera quality, date 44512
$ws.Cells.Item(126, 1).Value = 3
$ws.Cells.Item(126, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(126, 3).Value = 'Coquimbo'
$ws.Cells.Item(126, 4).Value = 44512
$ws.Cells.Item(126, 5).Value = 5
$ws.Cells.Item(126, 6).Value = 100112040
$ws.Cells.Item(126, 7).Value = 'Cilantro'
$ws.Cells.Item(126, 8).Value = 'Sin especificar'
$ws.Cells.Item(126, 9).Value = 'Primera'
$ws.Cells.Item(126, 10).Value = 130
$ws.Cells.Item(126, 11).Value = 3000
$ws.Cells.Item(126, 12).Value = 3000
$ws.Cells.Item(126, 13).Value = 3000
$ws.Cells.Item(126, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(126, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(126, 16).Value = 1000
$ws.Cells.Item(126, 17).Value = 3
$ws.Cells.Item(126, 18).Value = 'Hortaliza'

# Row 127 (new): Segunda quality, date 44512
$ws.Cells.Item(127, 1).Value = 3
$ws.Cells.Item(127, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(127, 3).Value = 'Coquimbo'
$ws.Cells.Item(127, 4).Value = 44512
$ws.Cells.Item(127, 5).Value = 5
$ws.Cells.Item(127, 6).Value = 100112040
$ws.Cells.Item(127, 7).Value = 'Cilantro'
$ws.Cells.Item(127, 8).Value = 'Sin especificar'
$ws.Cells.Item(127, 9).Value = 'Segunda'
$ws.Cells.Item(127, 10).Value = 90
$ws.Cells.Item(127, 11).Value = 2000
$ws.Cells.Item(127, 12).Value = 2000
$ws.Cells.Item(127, 13).Value = 2000
$ws.Cells.Item(127, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(127, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(127, 16).Value = 667
$ws.Cells.Item(127, 17).Value = 3
$ws.Cells.Item(127, 18).Value = 'Hortaliza'
